$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("B10").Value = "Fornecer ao aluno seminários sobre temas atuais de Física, Tecnologia e Engenharia."
$ws.Range("C10").Value = "Fornecer ao aluno seminários sobre temas atuais de Física, Tecnologia e Engenharia."

# Row 13
$ws.Range("A13").Clear()
$ws.Range("B13").Value = "3577649 - Carlos Angelo Nunes"
$ws.Range("C13").Value = "3577649 - Carlos Angelo Nunes"
$ws.Rows.Item(13).EntireRow.AutoFit()

# Row 14
$ws.Range("A14").Clear()
$ws.Range("B14").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C14").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Rows.Item(14).EntireRow.AutoFit()

# Row 15
$ws.Range("A15").Value = "Programa resumido:"
$ws.Range("B15").Value = "Seminários abrangendo os cenários atuais e futuros da indústria de alta tecnologia e do campo de atuação do engenheiro físico."
$ws.Range("C15").Value = "Seminários abrangendo os cenários atuais e futuros da indústria de alta tecnologia e do campo de atuação do engenheiro físico."
$ws.Rows.Item(15).RowHeight = 60

# Row 16
$ws.Range("A16").Value = "Short syllabus:"
$ws.Range("B16").Value = "Seminars covering the current and future scenarios of the high technology industry and the field of activity of the physical engineer."
$ws.Range("C16").Value = "Seminars covering the current and future scenarios of the high technology industry and the field of activity of the physical engineer."
$ws.Rows.Item(16).RowHeight = 60

# Row 17
$ws.Range("A17").Value = "Programa:"
$ws.Range("B17").Value = "Seminários seguido de debates com profissionais e estudantes de graduação e pós-graduação sobre temas relevantes e atuais das áreas de Física, Tecnologia e Engenharia, abrangendo desde as pesquisas básicas até o segmento industrial e de serviços."
$ws.Range("B3").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C17").Value = "Seminários seguido de debates com profissionais e estudantes de graduação e pós-graduação sobre temas relevantes e atuais das áreas de Física, Tecnologia e Engenharia, abrangendo desde as pesquisas básicas até o segmento industrial e de serviços."
$ws.Rows.Item(17).RowHeight = 120

# Row 18
$ws.Range("A18").Value = "Syllabus:"
$ws.Range("B18").Value = "Seminars followed by debates with professionals and undergraduate and graduate students on relevant and current topics in the areas of Physics, Technology and Engineering, ranging from basic research to the industrial and services segment."
$ws.Range("C18").Value = "Seminars followed by debates with professionals and undergraduate and graduate students on relevant and current topics in the areas of Physics, Technology and Engineering, ranging from basic research to the industrial and services segment."
$ws.Rows.Item(18).RowHeight = 120

# Row 19
$ws.Range("A19").Value = "Avaliação:"
$ws.Range("B19").Clear()
$ws.Range("C19").Clear()
$ws.Rows.Item(19).EntireRow.AutoFit()

# Row 20
$ws.Range("A20").Value = "Método:"
$ws.Range("B20").Value = "Os seminários proferidos por estudantes de graduação e pós-graduação, professores e convidados serão debatidos e analisados pelos alunos em forma de relatório. Os seminários apresentados pelos alunos serão avaliados na disciplina."
$ws.Range("C20").Value = "Os seminários proferidos por estudantes de graduação e pós-graduação, professores e convidados serão debatidos e analisados pelos alunos em forma de relatório. Os seminários apresentados pelos alunos serão avaliados na disciplina."
$ws.Rows.Item(20).EntireRow.AutoFit()

# Row 21
$ws.Range("A21").Value = "Critério:"
$ws.Range("B21").Value = "A nota final será calculada pela média aritmética dos relatórios e do seminário."
$ws.Range("C21").Value = "A nota final será calculada pela média aritmética dos relatórios e do seminário."
$ws.Rows.Item(21).RowHeight = 60

# Row 22
$ws.Range("A22").Value = "Norma de recuperação:"
$ws.Range("B22").Value = "Não há."
$ws.Range("B3").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C22").Value = "Não há."
$ws.Rows.Item(22).RowHeight = 60

# Row 23
$ws.Range("A23").Value = "Bibliografia:"
$ws.Range("B23").Value = "A ser definido de acordo com os temas dos seminários."
$ws.Range("C23").Value = "A ser definido de acordo com os temas dos seminários."
$ws.Rows.Item(23).RowHeight = 120

# Row 24
$ws.Range("A24").Value = "Requisitos:"

# Row 25
$ws.Range("B25").Value = "LOM3238 -  Projeto Integrado  (Requisito)`n"
$ws.Range("B3").Copy()
$ws.Range("B25").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C25").Value = "LOM3238 -  Projeto Integrado  (Requisito)`n"
$ws.Rows.Item(25).RowHeight = 30
